$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "<they>"
$ws.Range("C2").Value = 50

$ws.Range("B3").Value = "<her>"
$ws.Range("C3").Value = 49

$ws.Range("B4").Value = "<there>"
$ws.Range("C4").Value = 58

$ws.Range("B7").Value = "<oscar>"
$ws.Range("C7").Value = 52

$ws.Range("B8").Value = "<find>"
$ws.Range("C8").Value = 48

$ws.Range("B9").Value = "<her>"
$ws.Range("C9").Value = 49

$ws.Range("B10").Value = "<parte>"
$ws.Range("C10").Value = 56

$ws.Range("C11").Value = 52

$ws.Range("B12").Value = "<you>"
$ws.Range("C12").Value = 53

$ws.Range("B13").Value = "<on>"
$ws.Range("C13").Value = 57

$ws.Range("B14").Value = "<write>"
$ws.Range("C14").Value = 55

$ws.Range("B15").Value = "<four>"
$ws.Range("C15").Value = 52
